$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 31-61: corrected/refreshed interaction-score records and newly appended rows
$rows = @(
    @(31, "6738b019504ed0629a25b8b1", "67447fbced2b056beb0f8e01", 0.7, "2025-03-12T11:31:56.846Z"),
    @(32, "6738b019504ed0629a25b8b1", "676137906c06138b1419f8a5", 0.7, "2025-03-12T11:31:56.846Z"),
    @(33, "6738b019504ed0629a25b8b1", "67f08d651841d535b6af6f57", 0.7, "2025-03-12T11:31:56.846Z"),
    @(34, "6738b019504ed0629a25b8b1", "6728ebbb071b8fcf4f501e11", 0.7, "2025-03-12T11:31:56.846Z"),
    @(35, "6738b019504ed0629a25b8b1", "6728ea62071b8fcf4f501e02", 0.7, "2025-03-12T11:31:56.846Z"),
    @(36, "6738b019504ed0629a25b8b1", "6728e9ab071b8fcf4f501df6", 0.7, "2025-03-12T11:31:56.846Z"),
    @(37, "6738b019504ed0629a25b8b1", "67250625bb931ab886fc69db", 0.7, "2025-03-12T11:31:56.846Z"),
    @(38, "6738b019504ed0629a25b8b1", "6728e9e6071b8fcf4f501dfc", 0.7, "2025-03-12T11:31:56.846Z"),
    @(39, "6738b019504ed0629a25b8b1", "6728ea18071b8fcf4f501dff", 0.7, "2025-03-12T11:31:56.846Z"),
    @(40, "6738b019504ed0629a25b8b1", "6728eac6071b8fcf4f501e05", 0.7, "2025-03-12T11:31:56.846Z"),
    @(41, "6738b019504ed0629a25b8b1", "67f0959c1841d535b6af6fb4", 0.15, "2025-06-01T04:49:13.113Z"),
    @(42, "6738b019504ed0629a25b8b1", "67f094ec1841d535b6af6fab", 0.075, "2025-06-01T04:52:37.819Z"),
    @(43, "6738b019504ed0629a25b8b1", "6728e9cd071b8fcf4f501df9", 0.225, "2025-06-01T09:00:07.871Z"),
    @(44, "6738b019504ed0629a25b8b1", "67f095dc1841d535b6af6fba", 0.075, "2025-06-01T05:24:22.234Z"),
    @(45, "6738b019504ed0629a25b8b1", "67f094111841d535b6af6f99", 0.075, "2025-06-01T05:27:27.719Z"),
    @(46, "6738b019504ed0629a25b8b1", "67f092321841d535b6af6f81", 0.225, "2025-06-01T05:45:55.487Z"),
    @(47, "6738b019504ed0629a25b8b1", "6728e8a8071b8fcf4f501df0", 0.15, "2025-06-01T07:32:49.957Z"),
    @(48, "6738b019504ed0629a25b8b1", "67f095371841d535b6af6fb1", 0.3, "2025-06-01T07:48:03.296Z"),
    @(49, "6738b019504ed0629a25b8b1", "6728ec07071b8fcf4f501e17", 0.075, "2025-06-01T07:57:47.828Z"),
    @(50, "6738b019504ed0629a25b8b1", "6728ebe1071b8fcf4f501e14", 0.225, "2025-06-02T01:19:45.553Z"),
    @(51, "682326702fff19d415752f01", "6728e9e6071b8fcf4f501dfc", 0.925, "2025-05-23T04:41:14.096Z"),
    @(52, "682326702fff19d415752f01", "6728e93a071b8fcf4f501df3", 0.775, "2025-05-23T09:19:25.598Z"),
    @(53, "682326702fff19d415752f01", "6728ea18071b8fcf4f501dff", 0.775, "2025-05-23T09:19:36.004Z"),
    @(54, "682326702fff19d415752f01", "67f091181841d535b6af6f7b", 0.775, "2025-05-23T09:19:59.350Z"),
    @(55, "682326702fff19d415752f01", "6728e9cd071b8fcf4f501df9", 0.925, "2025-05-23T09:21:52.485Z"),
    @(56, "682326702fff19d415752f01", "6728f96acb86d3695fa1f4a6", 0.775, "2025-05-23T09:21:00.865Z"),
    @(57, "682326702fff19d415752f01", "68067dd1286f80e4174d8736", 1, "2025-05-23T09:23:22.245Z"),
    @(58, "683b37622eb85e2df9802771", "6728e93a071b8fcf4f501df3", 1, "2025-05-31T17:08:32.178Z"),
    @(59, "683b37622eb85e2df9802771", "6728ea18071b8fcf4f501dff", 0.85, "2025-05-31T17:10:43.229Z"),
    @(60, "683b37622eb85e2df9802771", "6728f96acb86d3695fa1f4a6", 0.15, "2025-06-01T02:53:50.928Z"),
    @(61, "683b37622eb85e2df9802771", "6728ec07071b8fcf4f501e17", 0.15, "2025-06-01T07:31:20.015Z")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}
